$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; rows 19-68 shift down to 20-69.
$ws.Rows("19:19").Insert()

# The new row 19 carries the same static/dimension columns as the row that
# used to be there (now row 20) but with new Fecha / Volumen / Precio data.
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = "2022-06-13"
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100108
$ws.Range("H19").Value = "Tropicales y subtropicales"
$ws.Range("I19").Value = 100108007
$ws.Range("J19").Value = "Coco"
$ws.Range("K19").Value = "Sin especificar"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 22000
$ws.Range("O19").Value = 23000
$ws.Range("P19").Value = 22500
$ws.Range("Q19").Value = "$/malla 20 unidades"
$ws.Range("R19").Value = "Perú"
$ws.Range("S19").Value = 1125
$ws.Range("T19").Value = 20
